$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Login")
$ws2 = $wb.Worksheets.Item("Cadastro")

# Update Login sheet: A2 value "Roger" -> "Victor"
$ws1.Range("A2").Value = "Victor"

# Update Cadastro sheet: A2 and C2 values
$ws2.Range("A2").Value = "AlbanoVictor17"
$ws2.Range("C2").Value = "AlbanoVictor17@gmail.com"

# Update selections/active cells (activate sheet first, then select the cell)
$ws2.Activate() | Out-Null
$ws2.Range("C2").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("B2").Select() | Out-Null
